$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 15 (shifts GBDT..XGBoost rows down by one)
$ws.Rows.Item(15).Insert()

# Set the new row's Model column to DeepCNN (other columns left blank)
$ws.Cells.Item(15, 1).Value = "DeepCNN"
